# "Model data now uses 4D vectors."
#
# The task "Move to 4D vectors in final model data" (row 2 - "A2:B2", the
# most recently added item at the top of the backlog) is removed from the
# ToDo list entirely: delete the whole row 2 and let every row below it
# slide up by one (values/styles move with their rows automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing review comment (anchored on the "Error handling
# strategy..." row, currently A8) so it can be re-anchored after the shift -
# row deletion does not automatically carry a cell comment to its new
# location.
$oldComment = $ws.Range("A8").Comment
$commentText = $oldComment.Text()
[void]$oldComment.Delete()

[void]$ws.Rows("2:2").Delete()

# Deleting a whole row leaves the selection spanning the full row that slid
# up into its place (matches Excel's own post-delete selection behaviour).
[void]$ws.Range("A2:XFD2").Select()

# Re-create the comment one row up (A8 -> A7) on the task it still belongs
# to ("Error handling strategy in Model Compiler - possibly use xsd?").
[void]$ws.Range("A7").AddComment($commentText)
